$d = $word.ActiveDocument
$x = $d.Styles.DefaultFont
if ($x -eq $null) { Write-Output "NULL" } else { Write-Output "NOTNULL: $($x.GetType())" }
